$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 6-16 (columns A: Player, B: Positions, C: Team)
$data = @(
    @("Paul George",           "SG,SF,PF", "Philadelphia 76ers"),
    @("Draymond Green",        "PF,C",     "Golden State Warriors"),
    @("Kyle Kuzma",            "PF",       "Washington Wizards"),
    @("Andrew Wiggins",        "SF,PF",    "Golden State Warriors"),
    @("Jaren Jackson Jr.",     "PF,C",     "Memphis Grizzlies"),
    @("Ivica Zubac",           "C",        "LA Clippers"),
    @("James Harden",          "PG,SG",    "LA Clippers"),
    @("Giannis Antetokounmpo", "PF,C",     "Milwaukee Bucks"),
    @("Caris LeVert",          "SG,SF",    "Cleveland Cavaliers"),
    @("Anthony Edwards",       "SG,SF",    "Minnesota Timberwolves"),
    @("Jayson Tatum",          "SF,PF",    "Boston Celtics")
)

$startRow = 6
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
